$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 is the header; row 4 (TetherUSD) is never touched by this update,
# so its cell keeps the plain default (un-styled) "Normal" style we can
# borrow from whenever we must force a numeric-looking string to stay text.
$normalStyle = $ws.Range("D4").Style

$ws.Range('D2').Value = '59.416.41'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '2.646.15'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '518.46'
$ws.Range('D5').Style = $normalStyle
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.25'
$ws.Range('D6').Style = $normalStyle
$ws.Range('E6').Value = '  -1.85%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = $normalStyle
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.572'
$ws.Range('D8').Style = $normalStyle
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('D9').Value = '2.650.84'
$ws.Range('E9').Value = '  -1.27%  '
$ws.Range('E10').Value = '  -3.61%  '
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.337'
$ws.Range('D12').Style = $normalStyle
$ws.Range('E12').Value = '  -1.49%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.128'
$ws.Range('D13').Style = $normalStyle
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').Value = '3.104.29'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').Value = '59.388.42'
$ws.Range('E15').Value = '  -0.69%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.95'
$ws.Range('D16').Style = $normalStyle
$ws.Range('E16').Value = '  -2.86%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000138'
$ws.Range('D17').Style = $normalStyle
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('D18').Value = '2.643.93'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '350.35'
$ws.Range('D19').Style = $normalStyle
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.50'
$ws.Range('D20').Style = $normalStyle
$ws.Range('E20').Value = '  -3.06%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.35'
$ws.Range('D21').Style = $normalStyle
$ws.Range('E21').Value = '  -3.02%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.21'
$ws.Range('D22').Style = $normalStyle
$ws.Range('E22').Value = '  -1.03%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = $normalStyle
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '62.03'
$ws.Range('D24').Style = $normalStyle
$ws.Range('E24').Value = '  +1.25%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.416'
$ws.Range('D25').Style = $normalStyle
$ws.Range('E25').Value = '  -3.06%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.165'
$ws.Range('D26').Style = $normalStyle
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.996'
$ws.Range('D27').Style = $normalStyle
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = '0.0₃0807'
$ws.Range('E28').Value = '  -3.67%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.14'
$ws.Range('D29').Style = $normalStyle
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').Style = $normalStyle
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.44'
$ws.Range('D31').Style = $normalStyle
$ws.Range('E31').Value = '  -2.92%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.58'
$ws.Range('D32').Style = $normalStyle
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.93'
$ws.Range('D33').Style = $normalStyle
$ws.Range('E33').Value = '  -1.06%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '149.27'
$ws.Range('D34').Style = $normalStyle
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.06'
$ws.Range('D35').Style = $normalStyle
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.940'
$ws.Range('D36').Style = $normalStyle
$ws.Range('E36').Value = '  -12.72%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.870'
$ws.Range('D38').Style = $normalStyle
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.61'
$ws.Range('D39').Style = $normalStyle
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.46'
$ws.Range('D40').Style = $normalStyle
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.67'
$ws.Range('D41').Style = $normalStyle
$ws.Range('E41').Value = '  -1.72%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '278.13'
$ws.Range('D42').Style = $normalStyle
$ws.Range('E42').Value = '  -4.28%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0992'
$ws.Range('D43').Style = $normalStyle
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.997'
$ws.Range('D44').Style = $normalStyle
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.602'
$ws.Range('D45').Style = $normalStyle
$ws.Range('E45').Value = '  -4.39%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '19.67'
$ws.Range('D46').Style = $normalStyle
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').Value = '2.106.50'
$ws.Range('E47').Value = '  +5.52%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0529'
$ws.Range('D48').Style = $normalStyle
$ws.Range('E48').Value = '  -4.06%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0232'
$ws.Range('D49').Style = $normalStyle
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.74'
$ws.Range('D50').Style = $normalStyle
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '10.30'
$ws.Range('D51').Style = $normalStyle
$ws.Range('E51').Value = '  -2.72%  '
